$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18 was a blank "<what day?>" template row. Turn it into a real diary
# entry (5th week - "mental simulation" lecture), following the same
# formatting pattern used by the other real entries.

# A18: date value, formatted/styled like the other date cells (A6/A7 use
# the "mm/dd/yyyy" + left-aligned style this row needs).
$ws.Range("A6").Copy()
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("A18").Value = 43867

# B18: time-range text, styled like D16:F16 (the other "font11" real rows).
$ws.Range("D16").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("B18").Value = "5:00pm - 8:00pm"

# C18 (Participants) stays blank for this entry - just drop the old
# placeholder text but keep the cell's existing look.
$ws.Range("C18").ClearContents()

# D18 (Goal), E18 (Achievements), F18 (Reflection): same style family as B18.
$ws.Range("D16").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("D18").Value = "Understand what was mental simulation.`nThe next 3 key expert practicies.`nHow do we use mental simulation with code"

$ws.Range("D16").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("E18").Value = "1) Learned what mental simulation was`n2) Learned that we use diagrams and concrete values as we step through code`n3) Learned to be skeptical of the code, of others, and of ourselves."

$ws.Range("D16").Copy()
$ws.Range("F18").PasteSpecial(-4122)
$ws.Range("F18").Value = "I learned that mental simulation was our way of imagining how the code works, and stepping through the code, thinking of the probably results without actually running the program. This was a way for us to tackle being skeptical of our code, since experts usually use print statements or debuggers to make sure that the value of something is actually that value. However, they simulate going through the code without running the program first. It is also helpful to draw diagrams and possible if-else situations, and stepping through it, simulating it. We have to ask ourselves what we are unsure of after every line of code."

# G18 (Your Overall Mood): styled like G16.
$ws.Range("G16").Copy()
$ws.Range("G18").PasteSpecial(-4122)
$ws.Range("G18").Value = "The activity gave another interesting perspective on how we deal with reading code, especially one within a large system that we have to understand. Like the mental models lecture, I felt that this lecture was also very helpful."
